$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, pushing all existing data down by one row
$ws.Rows.Item(1).Insert()

# Fill in the new header row
$ws.Range("A1").Value = "idquincena"
$ws.Range("B1").Value = "inicio"
$ws.Range("C1").Value = "fin"

# Apply right alignment to the new header row (matches style used for old header row)
$ws.Range("A1:C1").HorizontalAlignment = -4152

# Update selection to D8
$ws.Range("D8").Select()
